# "Quick Hack to include working hour update"
#
# Adds a new "Time Slot" column (column G) to the "driver" sheet:
#   - G1: header "Time Slot" (bold, matching the other header cells)
#   - G2: value 0
# and leaves the cursor/selection roughly where the author's saved file
# shows it (driver!G1:G2 selected while editing, then back to the
# vehicle sheet/tab with E1:E2 selected when the file was saved).

$wb = $excel.ActiveWorkbook

$driver = $wb.Worksheets.Item("driver")

# New "Time Slot" column header + value.
$driver.Range("G1").Value = "Time Slot"
$driver.Range("G1").Font.Bold = $true
$driver.Range("G2").Value = 0

# Leave the new column selected on the driver sheet ...
$driver.Range("G1:G2").Select()

# ... then hop back to "vehicle" (the tab that ends up active/saved)
# and leave its selection where the source file shows it.
$vehicle = $wb.Worksheets.Item("vehicle")
$vehicle.Activate()
$vehicle.Range("E1:E2").Select()
